# Regenerate Handback report: the localization pipeline ran again for the
# same two source files, producing new GUIDs/hashes and new timestamps.
# This mirrors the two file rows (7a435067... -> 291c5604..., and
# 9e8f5d6d... -> ffff58bf91f4...) across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper data: old -> new tokens
# ---------------------------------------------------------------------
# File 1: 7a435067-3689-4895-975e-8d657d9fe8af -> 291c5604-b16a-4a9d-9f37-928989d8f104
# File 2: 9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3 -> ffff58bf91f4-7e04-49b6-8eb0-f4826e2d7050
# Content hash: b90942244981dd365f980d378cfd591b010ac3f5 -> 47401c9d65c9b1103fca9549ff81d630b6a91ca0
#   (both rows now resolve to the SAME regenerated xlf content hash)

$file1Old = "7a435067-3689-4895-975e-8d657d9fe8af"
$file1New = "291c5604-b16a-4a9d-9f37-928989d8f104"
$file2Old = "9e8f5d6d-ca4a-467a-82d6-c76bcf3817f3"
$file2New = "ffff58bf91f4-7e04-49b6-8eb0-f4826e2d7050"
$hashNew  = "47401c9d65c9b1103fca9549ff81d630b6a91ca0"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$file1New.md"
$ws.Range("B2").Value = "e2e\$file1New.md"
$ws.Range("G2").Value = "2016-08-24 15:13:38"

$ws.Range("A3").Value = "$file2New.md"
$ws.Range("B3").Value = "e2e\$file2New.md"
$ws.Range("G3").Value = "2016-08-24 15:13:38"

# Refresh the hyperlink "display" text on B2/B3 to track the new file names
# (underlying link targets are left untouched, matching the source edit).
$b2Target = $ws.Hyperlinks.Item(1).Address
$b3Target = $ws.Hyperlinks.Item(2).Address

$ws.Range("B2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $b2Target, "", "", "e2e\$file1New.md")
$ws.Range("B2").Font.Underline = $true
$ws.Range("B2").Font.Color = 15570276

$ws.Range("B3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B3"), $b3Target, "", "", "e2e\$file2New.md")
$ws.Range("B3").Font.Underline = $true
$ws.Range("B3").Font.Color = 15570276

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$file1New.md"
$ws.Range("G2").Value = "$file1New.$hashNew.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-24 15:13:32"
$ws.Range("I2").Value = "$file1New.md"
$ws.Range("J2").Value = "$file1New.$hashNew.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-24 15:14:04"

$ws.Range("A3").Value = "$file2New.md"
$ws.Range("G3").Value = "$file1New.$hashNew.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-24 15:13:32"
$ws.Range("I3").Value = "$file2New.md"
$ws.Range("J3").Value = "$file1New.$hashNew.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-24 15:14:04"

$a2Target = $ws.Hyperlinks.Item(1).Address
$i2Target = $ws.Hyperlinks.Item(2).Address
$a3Target = $ws.Hyperlinks.Item(3).Address
$i3Target = $ws.Hyperlinks.Item(4).Address

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $a2Target, "", "", "$file1New.md")
$ws.Range("A2").Font.Underline = $true
$ws.Range("A2").Font.Color = 15570276

$ws.Range("I2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("I2"), $i2Target, "", "", "$file1New.md")
$ws.Range("I2").Font.Underline = $true
$ws.Range("I2").Font.Color = 15570276

$ws.Range("A3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A3"), $a3Target, "", "", "$file2New.md")
$ws.Range("A3").Font.Underline = $true
$ws.Range("A3").Font.Color = 15570276

$ws.Range("I3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("I3"), $i3Target, "", "", "$file2New.md")
$ws.Range("I3").Font.Underline = $true
$ws.Range("I3").Font.Color = 15570276

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$file1New.md"
$ws.Range("G2").Value = "$file1New.$hashNew.de-de.xlf"
$ws.Range("H2").Value = "2016-08-24 15:13:38"
$ws.Range("I2").Value = "$file1New.md"
$ws.Range("J2").Value = "$file1New.$hashNew.de-de.xlf"
$ws.Range("K2").Value = "2016-08-24 15:14:18"

$ws.Range("A3").Value = "$file2New.md"
$ws.Range("G3").Value = "$file1New.$hashNew.de-de.xlf"
$ws.Range("H3").Value = "2016-08-24 15:13:38"
$ws.Range("I3").Value = "$file2New.md"
$ws.Range("J3").Value = "$file1New.$hashNew.de-de.xlf"
$ws.Range("K3").Value = "2016-08-24 15:14:18"

$a2Target = $ws.Hyperlinks.Item(1).Address
$i2Target = $ws.Hyperlinks.Item(2).Address
$a3Target = $ws.Hyperlinks.Item(3).Address
$i3Target = $ws.Hyperlinks.Item(4).Address

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $a2Target, "", "", "$file1New.md")
$ws.Range("A2").Font.Underline = $true
$ws.Range("A2").Font.Color = 15570276

$ws.Range("I2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("I2"), $i2Target, "", "", "$file1New.md")
$ws.Range("I2").Font.Underline = $true
$ws.Range("I2").Font.Color = 15570276

$ws.Range("A3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A3"), $a3Target, "", "", "$file2New.md")
$ws.Range("A3").Font.Underline = $true
$ws.Range("A3").Font.Color = 15570276

$ws.Range("I3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("I3"), $i3Target, "", "", "$file2New.md")
$ws.Range("I3").Font.Underline = $true
$ws.Range("I3").Font.Color = 15570276
